$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1
$ws.Cells.Item(1, 2).Value = 45797
$ws.Cells.Item(1, 3).Value = 45804
$ws.Cells.Item(1, 4).Value = 45811
$ws.Cells.Item(1, 5).Value = 45818
$ws.Cells.Item(1, 6).Value = 45825
$ws.Cells.Item(1, 7).Value = 45832
$ws.Cells.Item(1, 8).Value = 45839
$ws.Cells.Item(1, 9).Value = 45846
$ws.Cells.Item(1, 10).Value = 45853
$ws.Cells.Item(1, 11).Value = 45860
$ws.Cells.Item(1, 12).Value = 45867

# Row 2
$ws.Cells.Item(2, 2).Value = 61.8
$ws.Cells.Item(2, 3).Value = 64.8
$ws.Cells.Item(2, 4).Value = 67
$ws.Cells.Item(2, 5).Value = 68.3
$ws.Cells.Item(2, 6).Value = 68.7
$ws.Cells.Item(2, 7).Value = 68.09999999999999
$ws.Cells.Item(2, 8).Value = 66.59999999999999
$ws.Cells.Item(2, 9).Value = 64.40000000000001
$ws.Cells.Item(2, 10).Value = 61.8
$ws.Cells.Item(2, 11).Value = 58.8
$ws.Cells.Item(2, 12).Value = 55.9

# Row 3
$ws.Cells.Item(3, 2).Value = 69.2
$ws.Cells.Item(3, 3).Value = 70.7
$ws.Cells.Item(3, 4).Value = 71.40000000000001
$ws.Cells.Item(3, 5).Value = 71.3
$ws.Cells.Item(3, 6).Value = 70
$ws.Cells.Item(3, 7).Value = 67.5
$ws.Cells.Item(3, 8).Value = 64
$ws.Cells.Item(3, 9).Value = 59.7
$ws.Cells.Item(3, 10).Value = 54.8
$ws.Cells.Item(3, 11).Value = 49.7
$ws.Cells.Item(3, 12).Value = 44.9

# Row 4
$ws.Cells.Item(4, 2).Value = 38.9
$ws.Cells.Item(4, 3).Value = 41.2
$ws.Cells.Item(4, 4).Value = 42.5
$ws.Cells.Item(4, 5).Value = 42.8
$ws.Cells.Item(4, 6).Value = 41.9
$ws.Cells.Item(4, 7).Value = 39.9
$ws.Cells.Item(4, 8).Value = 37
$ws.Cells.Item(4, 9).Value = 33.4
$ws.Cells.Item(4, 10).Value = 29.5
$ws.Cells.Item(4, 11).Value = 25.7
$ws.Cells.Item(4, 12).Value = 22.4

# Row 5
$ws.Cells.Item(5, 2).Value = 68.40000000000001
$ws.Cells.Item(5, 3).Value = 72
$ws.Cells.Item(5, 4).Value = 74.59999999999999
$ws.Cells.Item(5, 5).Value = 76.3
$ws.Cells.Item(5, 6).Value = 76.8
$ws.Cells.Item(5, 7).Value = 76.2
$ws.Cells.Item(5, 8).Value = 74.7
$ws.Cells.Item(5, 9).Value = 72.2
$ws.Cells.Item(5, 10).Value = 69.3
$ws.Cells.Item(5, 11).Value = 66
$ws.Cells.Item(5, 12).Value = 62.7

# Row 6
$ws.Cells.Item(6, 2).Value = 75.5
$ws.Cells.Item(6, 3).Value = 77.40000000000001
$ws.Cells.Item(6, 4).Value = 78.7
$ws.Cells.Item(6, 5).Value = 79.40000000000001
$ws.Cells.Item(6, 6).Value = 79.59999999999999
$ws.Cells.Item(6, 7).Value = 79.2
$ws.Cells.Item(6, 8).Value = 78.3
$ws.Cells.Item(6, 9).Value = 77.2
$ws.Cells.Item(6, 10).Value = 75.90000000000001
$ws.Cells.Item(6, 11).Value = 74.5
$ws.Cells.Item(6, 12).Value = 73.09999999999999

# Row 7
$ws.Cells.Item(7, 2).Value = 63.4
$ws.Cells.Item(7, 3).Value = 66.40000000000001
$ws.Cells.Item(7, 4).Value = 68.59999999999999
$ws.Cells.Item(7, 5).Value = 70
$ws.Cells.Item(7, 6).Value = 70.5
$ws.Cells.Item(7, 7).Value = 70.3
$ws.Cells.Item(7, 8).Value = 69.3
$ws.Cells.Item(7, 9).Value = 67.90000000000001
$ws.Cells.Item(7, 10).Value = 66.09999999999999
$ws.Cells.Item(7, 11).Value = 64.09999999999999
$ws.Cells.Item(7, 12).Value = 62

# Row 8
$ws.Cells.Item(8, 2).Value = 70.90000000000001
$ws.Cells.Item(8, 3).Value = 74.3
$ws.Cells.Item(8, 4).Value = 76.7
$ws.Cells.Item(8, 5).Value = 78
$ws.Cells.Item(8, 6).Value = 78.09999999999999
$ws.Cells.Item(8, 7).Value = 77
$ws.Cells.Item(8, 8).Value = 74.8
$ws.Cells.Item(8, 9).Value = 71.8
$ws.Cells.Item(8, 10).Value = 68.2
$ws.Cells.Item(8, 11).Value = 64.3
$ws.Cells.Item(8, 12).Value = 60.5

# Row 9
$ws.Cells.Item(9, 2).Value = 61.1
$ws.Cells.Item(9, 3).Value = 64.09999999999999
$ws.Cells.Item(9, 4).Value = 66.40000000000001
$ws.Cells.Item(9, 5).Value = 67.8
$ws.Cells.Item(9, 6).Value = 68.3
$ws.Cells.Item(9, 7).Value = 68
$ws.Cells.Item(9, 8).Value = 66.90000000000001
$ws.Cells.Item(9, 9).Value = 65.3
$ws.Cells.Item(9, 10).Value = 63.3
$ws.Cells.Item(9, 11).Value = 61.2
$ws.Cells.Item(9, 12).Value = 59.2

# Row 10
$ws.Cells.Item(10, 2).Value = 52.4
$ws.Cells.Item(10, 3).Value = 54.5
$ws.Cells.Item(10, 4).Value = 55.9
$ws.Cells.Item(10, 5).Value = 56.8
$ws.Cells.Item(10, 6).Value = 56.9
$ws.Cells.Item(10, 7).Value = 56.4
$ws.Cells.Item(10, 8).Value = 55.5
$ws.Cells.Item(10, 9).Value = 54.2
$ws.Cells.Item(10, 10).Value = 52.8
$ws.Cells.Item(10, 11).Value = 51.3
$ws.Cells.Item(10, 12).Value = 50

# Row 11
$ws.Cells.Item(11, 2).Value = 69.90000000000001
$ws.Cells.Item(11, 3).Value = 72.7
$ws.Cells.Item(11, 4).Value = 75.09999999999999
$ws.Cells.Item(11, 5).Value = 77
$ws.Cells.Item(11, 6).Value = 78.2
$ws.Cells.Item(11, 7).Value = 78.8
$ws.Cells.Item(11, 8).Value = 79
$ws.Cells.Item(11, 9).Value = 78.7
$ws.Cells.Item(11, 10).Value = 78.3
$ws.Cells.Item(11, 11).Value = 77.59999999999999
$ws.Cells.Item(11, 12).Value = 77

# Row 12
$ws.Cells.Item(12, 2).Value = 70
$ws.Cells.Item(12, 3).Value = 73.3
$ws.Cells.Item(12, 4).Value = 75.8
$ws.Cells.Item(12, 5).Value = 77.3
$ws.Cells.Item(12, 6).Value = 77.8
$ws.Cells.Item(12, 7).Value = 77.09999999999999
$ws.Cells.Item(12, 8).Value = 75.40000000000001
$ws.Cells.Item(12, 9).Value = 72.8
$ws.Cells.Item(12, 10).Value = 69.5
$ws.Cells.Item(12, 11).Value = 65.8
$ws.Cells.Item(12, 12).Value = 62.1

# Row 13
$ws.Cells.Item(13, 2).Value = 72.09999999999999
$ws.Cells.Item(13, 3).Value = 74.59999999999999
$ws.Cells.Item(13, 4).Value = 76.2
$ws.Cells.Item(13, 5).Value = 76.7
$ws.Cells.Item(13, 6).Value = 75.90000000000001
$ws.Cells.Item(13, 7).Value = 73.90000000000001
$ws.Cells.Item(13, 8).Value = 70.59999999999999
$ws.Cells.Item(13, 9).Value = 66.5
$ws.Cells.Item(13, 10).Value = 61.6
$ws.Cells.Item(13, 11).Value = 56.5
$ws.Cells.Item(13, 12).Value = 51.5

# Row 14
$ws.Cells.Item(14, 2).Value = 50.2
$ws.Cells.Item(14, 3).Value = 52.6
$ws.Cells.Item(14, 4).Value = 54.6
$ws.Cells.Item(14, 5).Value = 56
$ws.Cells.Item(14, 6).Value = 56.7
$ws.Cells.Item(14, 7).Value = 56.6
$ws.Cells.Item(14, 8).Value = 55.7
$ws.Cells.Item(14, 9).Value = 54.2
$ws.Cells.Item(14, 10).Value = 52.2
$ws.Cells.Item(14, 11).Value = 49.9
$ws.Cells.Item(14, 12).Value = 47.5

# Row 15
$ws.Cells.Item(15, 2).Value = 46.9
$ws.Cells.Item(15, 3).Value = 50.7
$ws.Cells.Item(15, 4).Value = 54.1
$ws.Cells.Item(15, 5).Value = 56.9
$ws.Cells.Item(15, 6).Value = 58.9
$ws.Cells.Item(15, 7).Value = 59.9
$ws.Cells.Item(15, 8).Value = 59.9
$ws.Cells.Item(15, 9).Value = 59
$ws.Cells.Item(15, 10).Value = 57.1
$ws.Cells.Item(15, 11).Value = 54.7
$ws.Cells.Item(15, 12).Value = 51.8

# Row 16
$ws.Cells.Item(16, 2).Value = 59.2
$ws.Cells.Item(16, 3).Value = 62.5
$ws.Cells.Item(16, 4).Value = 64.90000000000001
$ws.Cells.Item(16, 5).Value = 66.2
$ws.Cells.Item(16, 6).Value = 66.09999999999999
$ws.Cells.Item(16, 7).Value = 64.5
$ws.Cells.Item(16, 8).Value = 61.7
$ws.Cells.Item(16, 9).Value = 57.8
$ws.Cells.Item(16, 10).Value = 53.1
$ws.Cells.Item(16, 11).Value = 48.1
$ws.Cells.Item(16, 12).Value = 43.2

# Row 17
$ws.Cells.Item(17, 2).Value = 59.6
$ws.Cells.Item(17, 3).Value = 62.9
$ws.Cells.Item(17, 4).Value = 65.40000000000001
$ws.Cells.Item(17, 5).Value = 67
$ws.Cells.Item(17, 6).Value = 67.3
$ws.Cells.Item(17, 7).Value = 66.3
$ws.Cells.Item(17, 8).Value = 64.09999999999999
$ws.Cells.Item(17, 9).Value = 60.7
$ws.Cells.Item(17, 10).Value = 56.5
$ws.Cells.Item(17, 11).Value = 51.8
$ws.Cells.Item(17, 12).Value = 47

# Row 18
$ws.Cells.Item(18, 2).Value = 65.40000000000001
$ws.Cells.Item(18, 3).Value = 68.8
$ws.Cells.Item(18, 4).Value = 71.3
$ws.Cells.Item(18, 5).Value = 72.90000000000001
$ws.Cells.Item(18, 6).Value = 73.40000000000001
$ws.Cells.Item(18, 7).Value = 72.90000000000001
$ws.Cells.Item(18, 8).Value = 71.7
$ws.Cells.Item(18, 9).Value = 70
$ws.Cells.Item(18, 10).Value = 68
$ws.Cells.Item(18, 11).Value = 65.90000000000001
$ws.Cells.Item(18, 12).Value = 64

# Row 19
$ws.Cells.Item(19, 2).Value = 60.1
$ws.Cells.Item(19, 3).Value = 62.3
$ws.Cells.Item(19, 4).Value = 63.9
$ws.Cells.Item(19, 5).Value = 65.09999999999999
$ws.Cells.Item(19, 6).Value = 65.7
$ws.Cells.Item(19, 7).Value = 65.7
$ws.Cells.Item(19, 8).Value = 65.2
$ws.Cells.Item(19, 9).Value = 64.40000000000001
$ws.Cells.Item(19, 10).Value = 63.4
$ws.Cells.Item(19, 11).Value = 62.3
$ws.Cells.Item(19, 12).Value = 61.2

# Row 20
$ws.Cells.Item(20, 2).Value = 65.7
$ws.Cells.Item(20, 3).Value = 70.2
$ws.Cells.Item(20, 4).Value = 74.09999999999999
$ws.Cells.Item(20, 5).Value = 77.09999999999999
$ws.Cells.Item(20, 6).Value = 79.2
$ws.Cells.Item(20, 7).Value = 80.2
$ws.Cells.Item(20, 8).Value = 80.2
$ws.Cells.Item(20, 9).Value = 79.40000000000001
$ws.Cells.Item(20, 10).Value = 77.90000000000001
$ws.Cells.Item(20, 11).Value = 75.90000000000001
$ws.Cells.Item(20, 12).Value = 73.5

# Row 21
$ws.Cells.Item(21, 2).Value = 55.6
$ws.Cells.Item(21, 3).Value = 59.2
$ws.Cells.Item(21, 4).Value = 62.1
$ws.Cells.Item(21, 5).Value = 63.8
$ws.Cells.Item(21, 6).Value = 64.40000000000001
$ws.Cells.Item(21, 7).Value = 63.7
$ws.Cells.Item(21, 8).Value = 62
$ws.Cells.Item(21, 9).Value = 59.3
$ws.Cells.Item(21, 10).Value = 55.9
$ws.Cells.Item(21, 11).Value = 52.3
$ws.Cells.Item(21, 12).Value = 48.7
